$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 452, shifting existing rows 452:491 down to 453:492
$ws.Rows.Item(452).Insert()

# Populate the newly inserted row 452 with the new weekly record
$ws.Cells.Item(452, 1).Value = 7
$ws.Cells.Item(452, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(452, 3).Value = "Ñuble"
$ws.Cells.Item(452, 4).Value = 45132
$ws.Cells.Item(452, 5).Value = 16
$ws.Cells.Item(452, 6).Value = 100114013
$ws.Cells.Item(452, 7).Value = "Zanahoria"
$ws.Cells.Item(452, 8).Value = "Sin especificar"
$ws.Cells.Item(452, 9).Value = "Primera"
$ws.Cells.Item(452, 10).Value = 120
$ws.Cells.Item(452, 11).Value = 6000
$ws.Cells.Item(452, 12).Value = 6000
$ws.Cells.Item(452, 13).Value = 6000
$ws.Cells.Item(452, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(452, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(452, 16).Value = 300
$ws.Cells.Item(452, 17).Value = 20
$ws.Cells.Item(452, 18).Value = "Hortaliza"
